$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.321.38"
$ws.Range("E2").Value = "  -0.08%  "

$ws.Range("D3").Value = "3.686.96"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("D5").Value = "'681.01"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").Value = "'159.39"
$ws.Range("E6").Value = "  -1.98%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("E9").Value = "  -1.27%  "

$ws.Range("D10").Value = "'7.11"
$ws.Range("E10").Value = "  -3.29%  "

$ws.Range("D11").Value = "'0.439"
$ws.Range("E11").Value = "  -0.89%  "

$ws.Range("D12").Value = "'0.0000231"
$ws.Range("E12").Value = "  -3.34%  "

$ws.Range("D13").Value = "4.308.71"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("E14").Value = "  -2.68%  "

$ws.Range("D15").Value = "3.678.30"
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").Value = "69.314.51"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("E17").Value = "  +1.98%  "

$ws.Range("D18").Value = "'16.09"
$ws.Range("E18").Value = "  -0.90%  "

$ws.Range("D19").Value = "'6.46"
$ws.Range("E19").Value = "  -2.12%  "

$ws.Range("D20").Value = "'468.39"
$ws.Range("E20").Value = "  -2.42%  "

$ws.Range("D21").Value = "'9.94"
$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("E22").Value = "  -1.27%  "

$ws.Range("D23").Value = "'79.84"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("D24").Value = "3.833.27"
$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  -5.40%  "

$ws.Range("D27").Value = "'10.90"
$ws.Range("E27").Value = "  -4.04%  "

$ws.Range("E28").Value = "  -3.92%  "

$ws.Range("E29").Value = "  -1.39%  "

$ws.Range("E30").Value = "  -4.33%  "

$ws.Range("E31").Value = "  -2.82%  "

$ws.Range("E32").Value = "  -3.16%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("D35").Value = "3.675.82"
$ws.Range("E35").Value = "  +0.38%  "

$ws.Range("E36").Value = "  -5.18%  "

$ws.Range("E37").Value = "  -2.24%  "

$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  -2.72%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("E42").Value = "  -2.69%  "

$ws.Range("D43").Value = "'170.48"
$ws.Range("E43").Value = "  +3.90%  "

$ws.Range("D44").Value = "'0.943"
$ws.Range("E44").Value = "  -1.32%  "

$ws.Range("D45").Value = "'47.60"
$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("D46").Value = "'28.50"
$ws.Range("E46").Value = "  -5.52%  "

$ws.Range("E47").Value = "  -2.57%  "

$ws.Range("E48").Value = "  -4.01%  "

$ws.Range("D50").Value = "'0.000275"
$ws.Range("E50").Value = "  -3.58%  "

$ws.Range("E51").Value = "  -3.77%  "
